$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 4.286687
$ws.Range("H2").Value2 = 12.860061
$ws.Range("I2").Value2 = 0.3636010652499373
$ws.Range("J2").Value2 = 0.3636010652499373
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 1.316746666666667
$ws.Range("N2").Value2 = 3.95024
$ws.Range("O2").Value2 = 0.3073964349004645
$ws.Range("P2").Value2 = 0.3073964349004645
$ws.Range("Q2").Value2 = 5.644480818293333
$ws.Range("R2").Value2 = 50.80032736464
$ws.Range("S2").Value2 = 0.1117696711838419
$ws.Range("T2").Value2 = 0.1117696711838419

$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 4.286687
$ws.Range("H3").Value2 = 12.860061
$ws.Range("I3").Value2 = 0.3636010652499373
$ws.Range("J3").Value2 = 0.3636010652499373
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 1.660999
$ws.Range("N3").Value2 = 4.982997
$ws.Range("O3").Value2 = 0.3877626455404506
$ws.Range("P3").Value2 = 0.3877626455404507
$ws.Range("Q3").Value2 = 7.120182820313
$ws.Range("R3").Value2 = 64.081645382817
$ws.Range("S3").Value2 = 0.1409909109826417
$ws.Range("T3").Value2 = 0.1409909109826417

$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 4.286687
$ws.Range("H4").Value2 = 12.860061
$ws.Range("I4").Value2 = 0.3636010652499373
$ws.Range("J4").Value2 = 0.3636010652499373
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 1.3058
$ws.Range("N4").Value2 = 3.9174
$ws.Range("O4").Value2 = 0.3048409195590848
$ws.Range("P4").Value2 = 0.3048409195590849
$ws.Range("Q4").Value2 = 5.597555884599999
$ws.Range("R4").Value2 = 50.3780029614
$ws.Range("S4").Value2 = 0.1108404830834537
$ws.Range("T4").Value2 = 0.1108404830834537

$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 7.000795333333333
$ws.Range("H5").Value2 = 21.002386
$ws.Range("I5").Value2 = 0.5938144401018293
$ws.Range("J5").Value2 = 0.5938144401018293
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 1.316746666666667
$ws.Range("N5").Value2 = 3.95024
$ws.Range("O5").Value2 = 0.3073964349004645
$ws.Range("P5").Value2 = 0.3073964349004645
$ws.Range("Q5").Value2 = 9.218273919182224
$ws.Range("R5").Value2 = 82.96446527264
$ws.Range("S5").Value2 = 0.1825364418797177
$ws.Range("T5").Value2 = 0.1825364418797177

$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 7.000795333333333
$ws.Range("H6").Value2 = 21.002386
$ws.Range("I6").Value2 = 0.5938144401018293
$ws.Range("J6").Value2 = 0.5938144401018293
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 1.660999
$ws.Range("N6").Value2 = 4.982997
$ws.Range("O6").Value2 = 0.3877626455404506
$ws.Range("P6").Value2 = 0.3877626455404507
$ws.Range("Q6").Value2 = 11.62831404787133
$ws.Range("R6").Value2 = 104.654826430842
$ws.Range("S6").Value2 = 0.2302590582540068
$ws.Range("T6").Value2 = 0.2302590582540068

$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 7.000795333333333
$ws.Range("H7").Value2 = 21.002386
$ws.Range("I7").Value2 = 0.5938144401018293
$ws.Range("J7").Value2 = 0.5938144401018293
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 1.3058
$ws.Range("N7").Value2 = 3.9174
$ws.Range("O7").Value2 = 0.3048409195590848
$ws.Range("P7").Value2 = 0.3048409195590849
$ws.Range("Q7").Value2 = 9.141638546266666
$ws.Range("R7").Value2 = 82.2747469164
$ws.Range("S7").Value2 = 0.1810189399681047
$ws.Range("T7").Value2 = 0.1810189399681048

$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 0.5020513333333333
$ws.Range("H8").Value2 = 1.506154
$ws.Range("I8").Value2 = 0.04258449464823332
$ws.Range("J8").Value2 = 0.04258449464823332
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 1.316746666666667
$ws.Range("N8").Value2 = 3.95024
$ws.Range("O8").Value2 = 0.3073964349004645
$ws.Range("P8").Value2 = 0.3073964349004645
$ws.Range("Q8").Value2 = 0.6610744196622222
$ws.Range("R8").Value2 = 5.94966977696
$ws.Range("S8").Value2 = 0.01309032183690483
$ws.Range("T8").Value2 = 0.01309032183690483

$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 0.5020513333333333
$ws.Range("H9").Value2 = 1.506154
$ws.Range("I9").Value2 = 0.04258449464823332
$ws.Range("J9").Value2 = 0.04258449464823332
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 1.660999
$ws.Range("N9").Value2 = 4.982997
$ws.Range("O9").Value2 = 0.3877626455404506
$ws.Range("P9").Value2 = 0.3877626455404507
$ws.Range("Q9").Value2 = 0.8339067626153334
$ws.Range("R9").Value2 = 7.505160863538
$ws.Range("S9").Value2 = 0.01651267630380212
$ws.Range("T9").Value2 = 0.01651267630380212

$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 0.5020513333333333
$ws.Range("H10").Value2 = 1.506154
$ws.Range("I10").Value2 = 0.04258449464823332
$ws.Range("J10").Value2 = 0.04258449464823332
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 1.3058
$ws.Range("N10").Value2 = 3.9174
$ws.Range("O10").Value2 = 0.3048409195590848
$ws.Range("P10").Value2 = 0.3048409195590849
$ws.Range("Q10").Value2 = 0.6555786310666666
$ws.Range("R10").Value2 = 5.900207679599999
$ws.Range("S10").Value2 = 0.01298149650752637
$ws.Range("T10").Value2 = 0.01298149650752637

